$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,14
$row2[0,0] = 2.429188333333333
$row2[0,1] = 7.287565
$row2[0,2] = 0.1102134218696762
$row2[0,3] = 0.110840272037245
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 1.095307
$row2[0,7] = 3.285921
$row2[0,8] = 0.3507034061686541
$row2[0,9] = 0.3546600384155922
$row2[0,10] = 2.660706985818333
$row2[0,11] = 23.946362872365
$row2[0,12] = 0.03865222245519827
$row2[0,13] = 0.03931061513872401
$ws.Range("G2:T2").Value2 = $row2

$row3 = New-Object 'object[,]' 1,14
$row3[0,0] = 2.429188333333333
$row3[0,1] = 7.287565
$row3[0,2] = 0.1102134218696762
$row3[0,3] = 0.110840272037245
$row3[0,4] = 1
$row3[0,5] = 0.5
$row3[0,6] = 0.1045275
$row3[0,7] = 0.209055
$row3[0,8] = 0.03346837944822227
$row3[0,9] = 0.02256397957558067
$row3[0,10] = 0.2539169835125
$row3[0,11] = 1.523501901075
$row3[0,12] = 0.003688664623421321
$row3[0,13] = 0.002500997634400202
$ws.Range("G3:T3").Value2 = $row3

$row4 = New-Object 'object[,]' 1,14
$row4[0,0] = 2.429188333333333
$row4[0,1] = 7.287565
$row4[0,2] = 0.1102134218696762
$row4[0,3] = 0.110840272037245
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 1.923337333333333
$row4[0,7] = 5.770011999999999
$row4[0,8] = 0.6158282143831236
$row4[0,9] = 0.6227759820088272
$row4[0,10] = 4.672148611197778
$row4[0,11] = 42.04933750078
$row4[0,12] = 0.0678725347910566
$row4[0,13] = 0.06902865926412083
$ws.Range("G4:T4").Value2 = $row4

$row5 = New-Object 'object[,]' 1,14
$row5[0,0] = 2.839091333333334
$row5[0,1] = 8.517274
$row5[0,2] = 0.1288109145567312
$row5[0,3] = 0.1295435398759056
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 1.095307
$row5[0,7] = 3.285921
$row5[0,8] = 0.3507034061686541
$row5[0,9] = 0.3546600384155922
$row5[0,10] = 3.109676611039334
$row5[0,11] = 27.987089499354
$row5[0,12] = 0.04517442648674509
$row5[0,13] = 0.04594391682888049
$ws.Range("G5:T5").Value2 = $row5

$row6 = New-Object 'object[,]' 1,14
$row6[0,0] = 2.839091333333334
$row6[0,1] = 8.517274
$row6[0,2] = 0.1288109145567312
$row6[0,3] = 0.1295435398759056
$row6[0,4] = 1
$row6[0,5] = 0.5
$row6[0,6] = 0.1045275
$row6[0,7] = 0.209055
$row6[0,8] = 0.03346837944822227
$row6[0,9] = 0.02256397957558067
$row6[0,10] = 0.296763119345
$row6[0,11] = 1.78057871607
$row6[0,12] = 0.004311092565457217
$row6[0,13] = 0.002923017787908355
$ws.Range("G6:T6").Value2 = $row6

$row7 = New-Object 'object[,]' 1,14
$row7[0,0] = 2.839091333333334
$row7[0,1] = 8.517274
$row7[0,2] = 0.1288109145567312
$row7[0,3] = 0.1295435398759056
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 1.923337333333333
$row7[0,7] = 5.770011999999999
$row7[0,8] = 0.6158282143831236
$row7[0,9] = 0.6227759820088272
$row7[0,10] = 5.460530354143112
$row7[0,11] = 49.144773187288
$row7[0,12] = 0.07932539550452886
$row7[0,13] = 0.0806766052591168
$ws.Range("G7:T7").Value2 = $row7

$row8 = New-Object 'object[,]' 1,14
$row8[0,0] = 8.195700333333333
$row8[0,1] = 24.587101
$row8[0,2] = 0.3718427945500777
$row8[0,3] = 0.3739576886720351
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 1.095307
$row8[0,7] = 3.285921
$row8[0,8] = 0.3507034061686541
$row8[0,9] = 0.3546600384155922
$row8[0,10] = 8.976807945002333
$row8[0,11] = 80.791271505021
$row8[0,12] = 0.1304065346079833
$row8[0,13] = 0.13262784823023
$ws.Range("G8:T8").Value2 = $row8

$row9 = New-Object 'object[,]' 1,14
$row9[0,0] = 8.195700333333333
$row9[0,1] = 24.587101
$row9[0,2] = 0.3718427945500777
$row9[0,3] = 0.3739576886720351
$row9[0,4] = 1
$row9[0,5] = 0.5
$row9[0,6] = 0.1045275
$row9[0,7] = 0.209055
$row9[0,8] = 0.03346837944822227
$row9[0,9] = 0.02256397957558067
$row9[0,10] = 0.8566760665925
$row9[0,11] = 5.140056399555
$row9[0,12] = 0.01244497574308936
$row9[0,13] = 0.008437973649327156
$ws.Range("G9:T9").Value2 = $row9

$row10 = New-Object 'object[,]' 1,14
$row10[0,0] = 8.195700333333333
$row10[0,1] = 24.587101
$row10[0,2] = 0.3718427945500777
$row10[0,3] = 0.3739576886720351
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 1.923337333333333
$row10[0,7] = 5.770011999999999
$row10[0,8] = 0.6158282143831236
$row10[0,9] = 0.6227759820088272
$row10[0,10] = 15.76309642391244
$row10[0,11] = 141.867867815212
$row10[0,12] = 0.228991284199005
$row10[0,13] = 0.2328918667924779
$ws.Range("G10:T10").Value2 = $row10

$row11 = New-Object 'object[,]' 1,14
$row11[0,0] = 0.3739505
$row11[0,1] = 0.747901
$row11[0,2] = 0.01696631078345497
$row11[0,3] = 0.0113752056135249
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 1.095307
$row11[0,7] = 3.285921
$row11[0,8] = 0.3507034061686541
$row11[0,9] = 0.3546600384155922
$row11[0,10] = 0.4095906003035
$row11[0,11] = 2.457543601821
$row11[0,12] = 0.005950142981873625
$row11[0,13] = 0.004034330859878001
$ws.Range("G11:T11").Value2 = $row11

$row12 = New-Object 'object[,]' 1,14
$row12[0,0] = 0.3739505
$row12[0,1] = 0.747901
$row12[0,2] = 0.01696631078345497
$row12[0,3] = 0.0113752056135249
$row12[0,4] = 1
$row12[0,5] = 0.5
$row12[0,6] = 0.1045275
$row12[0,7] = 0.209055
$row12[0,8] = 0.03346837944822227
$row12[0,9] = 0.02256397957558067
$row12[0,10] = 0.03908811088875
$row12[0,11] = 0.156352443555
$row12[0,12] = 0.0005678349271371363
$row12[0,13] = 0.0002566699071316065
$ws.Range("G12:T12").Value2 = $row12

$row13 = New-Object 'object[,]' 1,14
$row13[0,0] = 0.3739505
$row13[0,1] = 0.747901
$row13[0,2] = 0.01696631078345497
$row13[0,3] = 0.0113752056135249
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 1.923337333333333
$row13[0,7] = 5.770011999999999
$row13[0,8] = 0.6158282143831236
$row13[0,9] = 0.6227759820088272
$row13[0,10] = 0.7192329574686667
$row13[0,11] = 4.315397744812
$row13[0,12] = 0.01044833287444421
$row13[0,13] = 0.007084204846515295
$ws.Range("G13:T13").Value2 = $row13

$row14 = New-Object 'object[,]' 1,14
$row14[0,0] = 8.202836333333334
$row14[0,1] = 24.608509
$row14[0,2] = 0.3721665582400601
$row14[0,3] = 0.3742832938012893
$row14[0,4] = 3
$row14[0,5] = 1
$row14[0,6] = 1.095307
$row14[0,7] = 3.285921
$row14[0,8] = 0.3507034061686541
$row14[0,9] = 0.3546600384155922
$row14[0,10] = 8.984624055754335
$row14[0,11] = 80.861616501789
$row14[0,12] = 0.1305200796368538
$row14[0,13] = 0.1327433273578796
$ws.Range("G14:T14").Value2 = $row14

$row15 = New-Object 'object[,]' 1,14
$row15[0,0] = 8.202836333333334
$row15[0,1] = 24.608509
$row15[0,2] = 0.3721665582400601
$row15[0,3] = 0.3742832938012893
$row15[0,4] = 1
$row15[0,5] = 0.5
$row15[0,6] = 0.1045275
$row15[0,7] = 0.209055
$row15[0,8] = 0.03346837944822227
$row15[0,9] = 0.02256397957558067
$row15[0,10] = 0.8574219748325
$row15[0,11] = 5.144531848995
$row15[0,12] = 0.01245581158911724
$row15[0,13] = 0.008445320596813351
$ws.Range("G15:T15").Value2 = $row15

$row16 = New-Object 'object[,]' 1,14
$row16[0,0] = 8.202836333333334
$row16[0,1] = 24.608509
$row16[0,2] = 0.3721665582400601
$row16[0,3] = 0.3742832938012893
$row16[0,4] = 3
$row16[0,5] = 1
$row16[0,6] = 1.923337333333333
$row16[0,7] = 5.770011999999999
$row16[0,8] = 0.6158282143831236
$row16[0,9] = 0.6227759820088272
$row16[0,10] = 15.77682135912311
$row16[0,11] = 141.991392232108
$row16[0,12] = 0.229190667014089
$row16[0,13] = 0.2330946458465964
$ws.Range("G16:T16").Value2 = $row16

